$p = $ppt.ActivePresentation

$oldDate = "11/6/2017"
$newDate = "11/7/2017"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder
Update-DateShapes($p.SlideMaster.Shapes)

# Every slide layout's date placeholder
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $cl = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DateShapes($cl.Shapes)
}

# Notes master date placeholder
Update-DateShapes($p.NotesMaster.Shapes)
